# Add a "Template" column (G) with VM image/template names used when
# creating each VM, to the Sizing/Config table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("G1").Value = "Template"

# Master / Worker / Storage / Bootstrap nodes all use the RHCOS image
$ws.Range("G2").Value = "rhcos-4.2.0-x86_64-vmware-template"
$ws.Range("G3").Value = "rhcos-4.2.0-x86_64-vmware-template"
$ws.Range("G4").Value = "rhcos-4.2.0-x86_64-vmware-template"
$ws.Range("G5").Value = "rhcos-4.2.0-x86_64-vmware-template"

# Install / LB / NFS nodes use their own dedicated templates
$ws.Range("G6").Value = "ocp42-installer-template"
$ws.Range("G7").Value = "ocp42-lb-template"
$ws.Range("G8").Value = "nfs-server-template"

# Widen the new column so the template names are readable
$ws.Columns("G").ColumnWidth = 34.833333333333336

# Update the sheet's current selection to cover the new column
[void]$ws.Range("A1:G8").Select()
